# Rectification de code de config
# Adds 9 new detail rows (142/143/146/147) before the TOTAUX row and
# updates the TOTAUX row accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: push the existing TOTAUX row (currently row 25) down
#        by 9 rows, so it ends up on row 34 and rows 25-33 become free
#        for the new detail lines. ---
$ws.Range("A25:A33").EntireRow.Insert()

# --- 2. Copy the formatting of the last existing detail row (24) onto
#        the freshly inserted rows so the new lines look like the rest
#        of the table (same borders / number formats as rows 2-24). ---
$ws.Range("A24:G24").Copy()
$ws.Range("A25:G33").PasteSpecial(-4122)

# --- 3. Fill in the new detail rows. Column A holds receipt numbers
#        that must stay text (leading apostrophe forces text storage
#        even though the value looks numeric). ---

# Row 25
$ws.Range("A25").Value = "'142"
$ws.Range("B25").Value = "MERC0001"
$ws.Range("C25").Value = "MERGUEZ CONGELE"
$ws.Range("D25").Value = 0.472
$ws.Range("E25").Value = 30000
$ws.Range("F25").Value = 14160
$ws.Range("G25").Value = "Espèces"

# Row 26
$ws.Range("A26").Value = "'142"
$ws.Range("B26").Value = "MERC0001"
$ws.Range("C26").Value = "MERGUEZ CONGELE"
$ws.Range("D26").Value = 0.39
$ws.Range("E26").Value = 30000
$ws.Range("F26").Value = 11700
$ws.Range("G26").Value = "Espèces"

# Row 27
$ws.Range("A27").Value = "'142"
$ws.Range("B27").Value = "MERL0001"
$ws.Range("C27").Value = " LAMB MERGUEZ / MERGUEZ D'AGNEAU "
$ws.Range("D27").Value = 0.33
$ws.Range("E27").Value = 60880
$ws.Range("F27").Value = 20090.4
$ws.Range("G27").Value = "Espèces"

# Row 28
$ws.Range("A28").Value = "'142"
$ws.Range("B28").Value = "SACHET_0"
$ws.Range("C28").Value = "PLASTIC BAGS"
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 500
$ws.Range("F28").Value = 500
$ws.Range("G28").Value = "Espèces"

# Row 29
$ws.Range("A29").Value = "'143"
$ws.Range("B29").Value = "BAF008"
$ws.Range("C29").Value = "RIBEYE / ENTRECOTE"
$ws.Range("D29").Value = 1.843
$ws.Range("E29").Value = 55510
$ws.Range("F29").Value = 102304.93
$ws.Range("G29").Value = "Espèces"

# Row 30
$ws.Range("A30").Value = "'143"
$ws.Range("B30").Value = "BAF008"
$ws.Range("C30").Value = "RIBEYE / ENTRECOTE"
$ws.Range("D30").Value = 0.547
$ws.Range("E30").Value = 55510
$ws.Range("F30").Value = 30363.97
$ws.Range("G30").Value = "Espèces"

# Row 31
$ws.Range("A31").Value = "'146"
$ws.Range("B31").Value = "LAML0001"
$ws.Range("C31").Value = "LAMB LEG CHOPS / GIGOT TRANCHE"
$ws.Range("D31").Value = 0.378
$ws.Range("E31").Value = 61950
$ws.Range("F31").Value = 23417.1
$ws.Range("G31").Value = "Espèces"

# Row 32
$ws.Range("A32").Value = "'146"
$ws.Range("B32").Value = "LAML0001"
$ws.Range("C32").Value = "LAMB LEG CHOPS / GIGOT TRANCHE"
$ws.Range("D32").Value = 0.31
$ws.Range("E32").Value = 61950
$ws.Range("F32").Value = 19204.5
$ws.Range("G32").Value = "Espèces"

# Row 33
$ws.Range("A33").Value = "'147"
$ws.Range("B33").Value = "FOICG001"
$ws.Range("C33").Value = "FOIE CONGELE"
$ws.Range("D33").Value = 0.292
$ws.Range("E33").Value = 17500
$ws.Range("F33").Value = 5110
$ws.Range("G33").Value = "Espèces"

# --- 4. Update the TOTAUX row, which is now row 34, with the new
#        totals that include the 9 additional lines. ---
$ws.Range("D34").Value = 15.034
$ws.Range("E34").Value = 1424467.5
$ws.Range("F34").Value = 634164.65
